# Loading time update after build success from jenkins
#
# Adds a new " Oct 21" column (X) to the loading-time report, mirroring the
# existing " Oct 19" column (W) values, with an updated reading for the
# "Deleted Trucks" row (row 6: 19 -> 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the workbook's plain "Normal" cell style before writing so the values
# land as real numbers/text instead of inheriting the sheet-wide Text number
# format (style 1) applied via the column default.
$ws.Range("X1:X7").Style = "Normal"

$ws.Range("X1").Value = " Oct 21"
$ws.Range("X2").Value = 0
$ws.Range("X3").Value = 5
$ws.Range("X4").Value = 5
$ws.Range("X5").Value = 0
$ws.Range("X6").Value = 22
$ws.Range("X7").Value = 0
